# Rename the workbook's first worksheet from "Sheet1" to "rates",
# matching the data that was imported (currency exchange rates) into
# the warehouse via SSIS.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "rates"
